$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new weekly record as row 180; all subsequent rows (old 180-214)
# shift down by one automatically (to 181-215), carrying their formatting/styles.
$ws.Rows.Item(180).Insert()

$newDate = Get-Date -Year 2023 -Month 3 -Day 30 -Hour 0 -Minute 0 -Second 0

$ws.Cells.Item(180, 1).Value = 8
$ws.Cells.Item(180, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(180, 3).Value = "Coquimbo"
$ws.Cells.Item(180, 4).Value = $newDate
$ws.Cells.Item(180, 5).Value = 4
$ws.Cells.Item(180, 6).Value = 100112044
$ws.Cells.Item(180, 7).Value = "Perejil"
$ws.Cells.Item(180, 8).Value = "Sin especificar"
$ws.Cells.Item(180, 9).Value = "Primera"
$ws.Cells.Item(180, 10).Value = 2400
$ws.Cells.Item(180, 11).Value = 1800
$ws.Cells.Item(180, 12).Value = 2000
$ws.Cells.Item(180, 13).Value = 1900
$ws.Cells.Item(180, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(180, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(180, 16).Value = 1267
$ws.Cells.Item(180, 17).Value = 1.5
$ws.Cells.Item(180, 18).Value = "Hortaliza"
